## Final Project Surveys.xlsx - "Add files via upload" edit
##
## The commit re-uploads the workbook after Excel filled in a bunch of
## previously-blank numeric cells with explicit 0 values (so every response
## row has a full row of zero-filled tallies instead of sparse blanks), and
## after a small "Extroversion" style table (A68:F78) picked up a thin-box
## border around every cell. The active sheet view also lost its scrolled
## position / old selection (back to the top-left, no stale selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

## ---------------------------------------------------------------
## 1) Explicit zero-fill for previously-empty cells.
##    (Cells that already hold a value are left completely alone.)
## ---------------------------------------------------------------

# Row 11-12 block (A:G of the 7-point Likert tally under header row 10)
$zeroCells = @(
    "B11",
    "C12", "F12",

    # Rows 15-24 (A:G tally block under header row 14)
    "B15", "C15", "D15", "E15", "F15", "G15",
    "B16", "D16", "E16", "F16", "G16",
    "B17", "D17", "F17", "G17",
    "E18", "F18", "G18",
    "E19",
    "G20",
    "C21", "D21", "G21",
    "D22",
    "B23", "C23", "G23",
    "B24", "C24", "E24",

    # Rows 26-28 (A:C tally block under header row 25)
    "B26", "C26",
    "B27", "C27",
    "C28",

    # Single category rows that only had a label, no count
    "B37",
    "B40",
    "B41",
    "B44",

    # Rows 69-78 (A:F "Extroversion" style block under header row 68)
    "B69", "C69", "D69", "E69", "F69",
    "B70", "C70", "E70", "F70",
    "D71", "F71",
    "C72", "D72", "F72",
    "F73",
    "D77", "E77",
    "B78"
)

foreach ($cellRef in $zeroCells) {
    $ws.Range($cellRef).Value = 0
}

## ---------------------------------------------------------------
## 2) Thin box border around the A68:F78 table (new style: borderId=1,
##    cellXfs applyBorder=1), matching the new <borders>/<cellXfs> entries.
## ---------------------------------------------------------------

$ws.Range("A68:F78").Borders.LineStyle = 1

## ---------------------------------------------------------------
## 3) Reset the sheet view: drop the scrolled-away topLeftCell and the
##    stale D6 selection, back to the natural top-left / A1.
## ---------------------------------------------------------------

$ws.Range("A1").Select()
